$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.5429463333333333
$ws.Cells.Item(2, 8).Value = 1.628839
$ws.Cells.Item(2, 9).Value = 0.04659251079363984
$ws.Cells.Item(2, 10).Value = 0.04659251079363985
$ws.Cells.Item(2, 13).Value = 2.231113333333334
$ws.Cells.Item(2, 14).Value = 6.69334
$ws.Cells.Item(2, 15).Value = 0.01598125358798882
$ws.Cells.Item(2, 16).Value = 0.01598125358798882
$ws.Cells.Item(2, 17).Value = 1.211374803584444
$ws.Cells.Item(2, 18).Value = 10.90237323226
$ws.Cells.Item(2, 19).Value = 0.0007446067302942645
$ws.Cells.Item(2, 20).Value = 0.0007446067302942646
$ws.Cells.Item(3, 7).Value = 0.5429463333333333
$ws.Cells.Item(3, 8).Value = 1.628839
$ws.Cells.Item(3, 9).Value = 0.04659251079363984
$ws.Cells.Item(3, 10).Value = 0.04659251079363985
$ws.Cells.Item(3, 15).Value = 0.1634493267640196
$ws.Cells.Item(3, 16).Value = 0.1634493267640195
$ws.Cells.Item(3, 17).Value = 12.38941582490033
$ws.Cells.Item(3, 18).Value = 111.504742424103
$ws.Cells.Item(3, 19).Value = 0.007615514521465748
$ws.Cells.Item(3, 20).Value = 0.007615514521465748
$ws.Cells.Item(4, 7).Value = 0.5429463333333333
$ws.Cells.Item(4, 8).Value = 1.628839
$ws.Cells.Item(4, 9).Value = 0.04659251079363984
$ws.Cells.Item(4, 10).Value = 0.04659251079363985
$ws.Cells.Item(4, 13).Value = 58.02175166666666
$ws.Cells.Item(4, 14).Value = 174.065255
$ws.Cells.Item(4, 15).Value = 0.4156043142904646
$ws.Cells.Item(4, 16).Value = 0.4156043142904646
$ws.Cells.Item(4, 17).Value = 31.50269732099389
$ws.Cells.Item(4, 18).Value = 283.5242758889449
$ws.Cells.Item(4, 19).Value = 0.01936404849946176
$ws.Cells.Item(4, 20).Value = 0.01936404849946176
$ws.Cells.Item(5, 7).Value = 0.5429463333333333
$ws.Cells.Item(5, 8).Value = 1.628839
$ws.Cells.Item(5, 9).Value = 0.04659251079363984
$ws.Cells.Item(5, 10).Value = 0.04659251079363985
$ws.Cells.Item(5, 13).Value = 15.16934033333333
$ws.Cells.Item(5, 14).Value = 45.508021
$ws.Cells.Item(5, 15).Value = 0.1086565487318021
$ws.Cells.Item(5, 16).Value = 0.1086565487318021
$ws.Cells.Item(5, 17).Value = 8.236137713068777
$ws.Cells.Item(5, 18).Value = 74.125239417619
$ws.Cells.Item(5, 19).Value = 0.005062581419586145
$ws.Cells.Item(5, 20).Value = 0.005062581419586145
$ws.Cells.Item(6, 7).Value = 0.5429463333333333
$ws.Cells.Item(6, 8).Value = 1.628839
$ws.Cells.Item(6, 9).Value = 0.04659251079363984
$ws.Cells.Item(6, 10).Value = 0.04659251079363985
$ws.Cells.Item(6, 13).Value = 41.36709099999999
$ws.Cells.Item(6, 14).Value = 124.101273
$ws.Cells.Item(6, 15).Value = 0.2963085566257249
$ws.Cells.Item(6, 16).Value = 0.2963085566257249
$ws.Cells.Item(6, 17).Value = 22.46011037911633
$ws.Cells.Item(6, 18).Value = 202.140993412047
$ws.Cells.Item(6, 19).Value = 0.01380575962283193
$ws.Cells.Item(6, 20).Value = 0.01380575962283193
$ws.Cells.Item(7, 9).Value = 0.8858267105024722
$ws.Cells.Item(7, 10).Value = 0.8858267105024723
$ws.Cells.Item(7, 13).Value = 2.231113333333334
$ws.Cells.Item(7, 14).Value = 6.69334
$ws.Cells.Item(7, 15).Value = 0.01598125358798882
$ws.Cells.Item(7, 16).Value = 0.01598125358798882
$ws.Cells.Item(7, 17).Value = 23.03091503691333
$ws.Cells.Item(7, 18).Value = 207.27823533222
$ws.Cells.Item(7, 19).Value = 0.01415662129555397
$ws.Cells.Item(7, 20).Value = 0.01415662129555397
$ws.Cells.Item(8, 9).Value = 0.8858267105024722
$ws.Cells.Item(8, 10).Value = 0.8858267105024723
$ws.Cells.Item(8, 15).Value = 0.1634493267640196
$ws.Cells.Item(8, 16).Value = 0.1634493267640195
$ws.Cells.Item(8, 19).Value = 0.1447877794612152
$ws.Cells.Item(8, 20).Value = 0.1447877794612152
$ws.Cells.Item(9, 9).Value = 0.8858267105024722
$ws.Cells.Item(9, 10).Value = 0.8858267105024723
$ws.Cells.Item(9, 13).Value = 58.02175166666666
$ws.Cells.Item(9, 14).Value = 174.065255
$ws.Cells.Item(9, 15).Value = 0.4156043142904646
$ws.Cells.Item(9, 16).Value = 0.4156043142904646
$ws.Cells.Item(9, 17).Value = 598.9359719936016
$ws.Cells.Item(9, 18).Value = 5390.423747942415
$ws.Cells.Item(9, 19).Value = 0.3681534025985578
$ws.Cells.Item(9, 20).Value = 0.3681534025985579
$ws.Cells.Item(10, 9).Value = 0.8858267105024722
$ws.Cells.Item(10, 10).Value = 0.8858267105024723
$ws.Cells.Item(10, 13).Value = 15.16934033333333
$ws.Cells.Item(10, 14).Value = 45.508021
$ws.Cells.Item(10, 15).Value = 0.1086565487318021
$ws.Cells.Item(10, 16).Value = 0.1086565487318021
$ws.Cells.Item(10, 17).Value = 156.5871993876103
$ws.Cells.Item(10, 18).Value = 1409.284794488493
$ws.Cells.Item(10, 19).Value = 0.09625087313764384
$ws.Cells.Item(10, 20).Value = 0.09625087313764386
$ws.Cells.Item(11, 9).Value = 0.8858267105024722
$ws.Cells.Item(11, 10).Value = 0.8858267105024723
$ws.Cells.Item(11, 13).Value = 41.36709099999999
$ws.Cells.Item(11, 14).Value = 124.101273
$ws.Cells.Item(11, 15).Value = 0.2963085566257249
$ws.Cells.Item(11, 16).Value = 0.2963085566257249
$ws.Cells.Item(11, 17).Value = 427.016388594601
$ws.Cells.Item(11, 18).Value = 3843.147497351409
$ws.Cells.Item(11, 19).Value = 0.2624780340095014
$ws.Cells.Item(11, 20).Value = 0.2624780340095014
$ws.Cells.Item(12, 7).Value = 0.7875243333333334
$ws.Cells.Item(12, 8).Value = 2.362573
$ws.Cells.Item(12, 9).Value = 0.06758077870388791
$ws.Cells.Item(12, 10).Value = 0.06758077870388793
$ws.Cells.Item(12, 13).Value = 2.231113333333334
$ws.Cells.Item(12, 14).Value = 6.69334
$ws.Cells.Item(12, 15).Value = 0.01598125358798882
$ws.Cells.Item(12, 16).Value = 0.01598125358798882
$ws.Cells.Item(12, 17).Value = 1.757056040424445
$ws.Cells.Item(12, 18).Value = 15.81350436382
$ws.Cells.Item(12, 19).Value = 0.001080025562140587
$ws.Cells.Item(12, 20).Value = 0.001080025562140587
$ws.Cells.Item(13, 7).Value = 0.7875243333333334
$ws.Cells.Item(13, 8).Value = 2.362573
$ws.Cells.Item(13, 9).Value = 0.06758077870388791
$ws.Cells.Item(13, 10).Value = 0.06758077870388793
$ws.Cells.Item(13, 15).Value = 0.1634493267640196
$ws.Cells.Item(13, 16).Value = 0.1634493267640195
$ws.Cells.Item(13, 17).Value = 17.97040672140233
$ws.Cells.Item(13, 18).Value = 161.733660492621
$ws.Cells.Item(13, 19).Value = 0.01104603278133867
$ws.Cells.Item(13, 20).Value = 0.01104603278133867
$ws.Cells.Item(14, 7).Value = 0.7875243333333334
$ws.Cells.Item(14, 8).Value = 2.362573
$ws.Cells.Item(14, 9).Value = 0.06758077870388791
$ws.Cells.Item(14, 10).Value = 0.06758077870388793
$ws.Cells.Item(14, 13).Value = 58.02175166666666
$ws.Cells.Item(14, 14).Value = 174.065255
$ws.Cells.Item(14, 15).Value = 0.4156043142904646
$ws.Cells.Item(14, 16).Value = 0.4156043142904646
$ws.Cells.Item(14, 17).Value = 45.69354130012388
$ws.Cells.Item(14, 18).Value = 411.241871701115
$ws.Cells.Item(14, 19).Value = 0.02808686319244497
$ws.Cells.Item(14, 20).Value = 0.02808686319244497
$ws.Cells.Item(15, 7).Value = 0.7875243333333334
$ws.Cells.Item(15, 8).Value = 2.362573
$ws.Cells.Item(15, 9).Value = 0.06758077870388791
$ws.Cells.Item(15, 10).Value = 0.06758077870388793
$ws.Cells.Item(15, 13).Value = 15.16934033333333
$ws.Cells.Item(15, 14).Value = 45.508021
$ws.Cells.Item(15, 15).Value = 0.1086565487318021
$ws.Cells.Item(15, 16).Value = 0.1086565487318021
$ws.Cells.Item(15, 17).Value = 11.94622463311478
$ws.Cells.Item(15, 18).Value = 107.516021698033
$ws.Cells.Item(15, 19).Value = 0.007343094174572132
$ws.Cells.Item(15, 20).Value = 0.007343094174572134
$ws.Cells.Item(16, 7).Value = 0.7875243333333334
$ws.Cells.Item(16, 8).Value = 2.362573
$ws.Cells.Item(16, 9).Value = 0.06758077870388791
$ws.Cells.Item(16, 10).Value = 0.06758077870388793
$ws.Cells.Item(16, 13).Value = 41.36709099999999
$ws.Cells.Item(16, 14).Value = 124.101273
$ws.Cells.Item(16, 15).Value = 0.2963085566257249
$ws.Cells.Item(16, 16).Value = 0.2963085566257249
$ws.Cells.Item(16, 17).Value = 32.57759076171433
$ws.Cells.Item(16, 18).Value = 293.198316855429
$ws.Cells.Item(16, 19).Value = 0.02002476299339156
$ws.Cells.Item(16, 20).Value = 0.02002476299339156
